$d = $word.ActiveDocument

# The edit adds a new paragraph right after the first paragraph
# ("Cognitive performance: add ... after CN_ADNI"), containing the
# holdout CU correlation results. The existing "_GoBack" bookmark that
# sits at the end of paragraph 1 needs to move to the end of the new
# paragraph, so we delete it first and re-create it inside the XML we
# insert below.
$p1 = $d.Paragraphs(1)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Collapsed range right before the paragraph mark that ends paragraph 1.
$splitPos = $p1.Range.End - 1
$insertionPoint = $d.Range($splitPos, $splitPos)

$newParagraphXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:vertAlign w:val="subscript"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t>CU</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t>: n</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t>=</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">237, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t>ADNI-EF</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t>r</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:vertAlign w:val="subscript"/><w:lang w:val="en-US"/></w:rPr><w:t>MRI</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:vertAlign w:val="subscript"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">= .023, p = .73, 95% CI [-.11, .15]; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t>r</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:vertAlign w:val="subscript"/><w:lang w:val="en-US"/></w:rPr><w:t>FDG</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:vertAlign w:val="subscript"/><w:lang w:val="en-US"/></w:rPr><w:t>-PET</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> = -.019, p = .77, 95% CI [-.15, .11]; </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t>ADNI-MEM</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t>r</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:vertAlign w:val="subscript"/><w:lang w:val="en-US"/></w:rPr><w:t>MRI</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:vertAlign w:val="subscript"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">= -.048, p = .46, 95% CI [-.18, .08]; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t>r</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:vertAlign w:val="subscript"/><w:lang w:val="en-US"/></w:rPr><w:t>FDG</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:vertAlign w:val="subscript"/><w:lang w:val="en-US"/></w:rPr><w:t>-PET</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> = -.015, p = .82, 95% CI [-.14</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t>,0.11</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t>]</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$insertionPoint.InsertXML($newParagraphXml) | Out-Null
